$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price/volume data (GitHub Actions scheduled update).
# Column D ("Price") values that look numeric must be forced to Text so
# Excel keeps literal strings like "1.00" / "0.999" instead of collapsing
# them to plain numbers, matching the source data formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.577.78"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.516.23"
$ws.Range("E3").Value = "  -2.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.93"
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.83"
$ws.Range("E6").Value = "  -4.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.515.74"
$ws.Range("E7").Value = "  -2.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("E9").Value = "  +3.95%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.74"
$ws.Range("E10").Value = "  -3.79%  "

$ws.Range("E11").Value = "  -4.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.407"
$ws.Range("E12").Value = "  -2.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.112.57"
$ws.Range("E13").Value = "  -2.29%  "

$ws.Range("E14").Value = "  -6.78%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.59"
$ws.Range("E15").Value = "  -4.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.512.12"
$ws.Range("E16").Value = "  -2.62%  "

$ws.Range("E17").Value = "  +0.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.416.28"
$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.77"
$ws.Range("E19").Value = "  -6.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.15"
$ws.Range("E20").Value = "  -3.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.68"
$ws.Range("E21").Value = "  -2.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "421.55"
$ws.Range("E22").Value = "  -1.33%  "

$ws.Range("E23").Value = "  -5.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "76.68"
$ws.Range("E24").Value = "  -2.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.655.42"
$ws.Range("E25").Value = "  -2.28%  "

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("E27").Value = "  -6.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.92"
$ws.Range("E28").Value = "  -4.90%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.46"
$ws.Range("E29").Value = "  -2.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.91"
$ws.Range("E30").Value = "  -5.72%  "

$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.521.96"
$ws.Range("E32").Value = "  -2.06%  "

$ws.Range("E33").Value = "  -2.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.16"
$ws.Range("E34").Value = "  -5.08%  "

$ws.Range("E36").Value = "  -9.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.55"
$ws.Range("E37").Value = "  -3.94%  "

$ws.Range("E38").Value = "  -4.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "173.42"
$ws.Range("E39").Value = "  -2.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.21"
$ws.Range("E40").Value = "  -7.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0816"
$ws.Range("E41").Value = "  -4.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.97"
$ws.Range("E42").Value = "  -5.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.852"
$ws.Range("E43").Value = "  -5.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.48"
$ws.Range("E44").Value = "  -0.90%  "

$ws.Range("E45").Value = "  -7.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("E47").Value = "  -8.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.05"
$ws.Range("E48").Value = "  -2.16%  "

$ws.Range("E49").Value = "  -5.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.74"
$ws.Range("E50").Value = "  -5.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.903"
$ws.Range("E51").Value = "  -5.26%  "
